$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their text formatting so values are stored as strings
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '307.34'
$ws.Range('E2').Value = '-0.81%'
$ws.Range('D3').Value = '36.95'
$ws.Range('E3').Value = '-0.93%'
$ws.Range('D4').Value = '5.113'
$ws.Range('E4').Value = '-0.25%'
$ws.Range('D5').Value = '0.07804'
$ws.Range('E5').Value = '0.47%'
$ws.Range('D6').Value = '8.208'
$ws.Range('E6').Value = '-0.07%'
$ws.Range('D7').Value = '1.879'
$ws.Range('E7').Value = '-0.09%'
$ws.Range('E8').Value = '-2.18%'
$ws.Range('D9').Value = '0.9209'
$ws.Range('E9').Value = '0.25%'
$ws.Range('D10').Value = '0.1087'
$ws.Range('E10').Value = '-9.03%'
$ws.Range('D11').Value = '0.1894'
$ws.Range('E11').Value = '-0.15%'
$ws.Range('D12').Value = '0.08900'
$ws.Range('E12').Value = '-2.81%'
$ws.Range('D13').Value = '0.03348'
$ws.Range('E13').Value = '-2.21%'
$ws.Range('D14').Value = '0.09574'
$ws.Range('E14').Value = '-1.20%'
$ws.Range('D15').Value = '0.001380'
$ws.Range('E15').Value = '0.94%'
$ws.Range('D16').Value = '0.005804'
$ws.Range('E16').Value = '-1.12%'
$ws.Range('D17').Value = '3.442'
$ws.Range('E17').Value = '-3.12%'
$ws.Range('D18').Value = '4.393'
$ws.Range('E18').Value = '0.11%'
$ws.Range('E19').Value = '0.50%'
$ws.Range('D20').Value = '6.227'
$ws.Range('E20').Value = '18.47%'
$ws.Range('D21').Value = '0.1286'
$ws.Range('E21').Value = '0.90%'
$ws.Range('D22').Value = '0.2433'
$ws.Range('E22').Value = '-6.19%'
$ws.Range('D23').Value = '0.04345'
$ws.Range('E23').Value = '-0.48%'
$ws.Range('D24').Value = '0.001194'
$ws.Range('E24').Value = '-0.45%'
$ws.Range('D25').Value = '0.004255'
$ws.Range('E25').Value = '0.14%'
$ws.Range('E26').Value = '0.67%'
$ws.Range('E27').Value = '-98.10%'
$ws.Range('D39').Value = '0.02157'
$ws.Range('E39').Value = '3.46%'
$ws.Range('D40').Value = '0.05016'
$ws.Range('E40').Value = '-0.11%'
$ws.Range('D41').Value = '0.007525'
$ws.Range('E41').Value = '-1.53%'
$ws.Range('D42').Value = '0.1347'
$ws.Range('D43').Value = '0.008672'
$ws.Range('E43').Value = '-11.97%'
$ws.Range('D44').Value = '0.002031'
$ws.Range('E44').Value = '-6.54%'
$ws.Range('D45').Value = '0.008737'
$ws.Range('E45').Value = '-9.00%'
$ws.Range('D46').Value = '0.00006534'
$ws.Range('E46').Value = '-2.63%'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').Value = '-0.07%'
$ws.Range('D48').Value = '0.003379'
$ws.Range('E48').Value = '15.03%'
$ws.Range('E49').Value = '-16.59%'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('E50').Value = '-0.07%'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('E51').Value = '-0.07%'
